$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Insert two new paragraphs ("Abstract: " and the justified abstract
# body) followed by 28 blank paragraphs, all ahead of the existing
# "Draft of the FOOOF paper " paragraph that used to open the body.
# ------------------------------------------------------------------

# Make room: two blank paragraphs land at the very top of the body.
$d.Range(0, 0).InsertParagraphBefore()
$d.Range(0, 0).InsertParagraphBefore()

# Paragraph 1: "Abstract: "
$d.Paragraphs(1).Range.Text = "Abstract: "

# Paragraph 2: the justified abstract body, rebuilt run by run so the
# text matches the source edit exactly.
$abstractRuns = @(
    "There is a significant push for using parametric description of Neural ",
    "Power ",
    "S",
    "pectrum",
    ". ",
    "In particular, one parametric model, FOOOF, has been ",
    "introduced recently and it gained significant attention in the field. ",
    "FOOOF algorithm suggested a parametric model to estimate ",
    "Power Spectrum density of signals. However, ",
    "the FOOOF algorithm provides an interesting framework to assess ",
    "frequency domain ",
    "time-series data",
    "; it suffers from numerous of shortages. ",
    "FOOOF model is not able to provide a robust frequency estimation, ",
    "control variability of parameters, and",
    " include the notion of continuity of neural signals in the model. ",
    "In this research, we are ",
    "introducing",
    " a Bayesian ",
    "FOOOF model which addresses multiple issues of the previous model including time continuity over time, ",
    "much more flexibility in controlling ",
    "specific frequency bands, and ",
    "also rather",
    " than providing a point estimate for each free parameter",
    ", this model will provide an posterior estimation of each parameter of the model. ",
    "In addition to the method development information of this proposed model, we also ",
    "develop a toolbox which can be used for lots of different kinds of time series data without having the ",
    "expertise of the field. This model would be a significant endeavor in the ",
    "computational neuroscience field",
    " to provide a parametric model ",
    "of the Power Spectrum Density",
    ", which is very important notion in Neuroscience data analysis. "
)

$pos = $d.Paragraphs(2).Range.Start
foreach ($runText in $abstractRuns) {
    $r = $d.Range($pos, $pos)
    $r.InsertBefore($runText)
    $pos = $pos + $runText.Length
}
$d.Paragraphs(2).Alignment = 3  # wdAlignParagraphJustify -> <w:jc w:val="both"/>

# 28 blank paragraphs between the abstract and the old first paragraph.
$oldFirstParagraph = $d.Paragraphs(3).Range
for ($i = 0; $i -lt 28; $i++) {
    $oldFirstParagraph.InsertParagraphBefore()
}

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
